$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new header "REX_DEF" in F1, matching the formatting of the other
# header cells (B1:E1) by copying E1's format onto F1.
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("F1").Value = "REX_DEF"

# Fill F2:F10 with the literal value "[]" (unstyled, like column C/E data cells)
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 6).Value = "[]"
}
